# Houston roster update: swap three pairs of player rows so the roster
# table reflects the corrected ordering (players swapped places):
#   Row 5  <-> Row 6   (Alperen Sengun  <-> Jalen Green)
#   Row 9  <-> Row 10  (Daishen Nix     <-> Kevin Porter Jr.)
#   Row 11 <-> Row 12  (TyTy Washington Jr. <-> Jae'Sean Tate)
#
# Columns in the table:
#   A=index  B=No.  C=Player  D=Pos  E=Ht  F=Wt  G=Birth Date
#   H=nationality code  I=Exp  J=College  K=bbref url
# (Column A, the running 0-based index, is NOT part of the swap - it stays
# attached to the row.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value2
        $valB = $rangeB.Value2
        $rangeA.Value = $valB
        $rangeB.Value = $valA
    }
}

Swap-Rows 5 6
Swap-Rows 9 10
Swap-Rows 11 12
